$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5334
$ws.Range("F6").Value = 221
$ws.Range("F7").Value = 186
$ws.Range("F8").Value = 8939
$ws.Range("G8").Value = 238
$ws.Range("F9").Value = 278
$ws.Range("F11").Value = 13
$ws.Range("F12").Value = 2621
$ws.Range("F13").Value = 2621
$ws.Range("F14").Value = 6348
$ws.Range("F15").Value = 2347
$ws.Range("F19").Value = 2549
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 6613
$ws.Range("F23").Value = 223
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 158
$ws.Range("F28").Value = 7267
$ws.Range("F31").Value = 244
$ws.Range("F32").Value = 45
$ws.Range("F35").Value = 29
$ws.Range("F40").Value = 2555
$ws.Range("F43").Value = 16
$ws.Range("F46").Value = 560
$ws.Range("F47").Value = 3574
$ws.Range("F48").Value = 107
$ws.Range("F49").Value = 1142
$ws.Range("F50").Value = 42

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 26
$ws.Range("F5").Value = 216
$ws.Range("F7").Value = 100

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5336
$ws.Range("F4").Value = 5336
$ws.Range("F6").Value = 221
$ws.Range("F7").Value = 8939
$ws.Range("G7").Value = 238
$ws.Range("F8").Value = 278
$ws.Range("F10").Value = 26
$ws.Range("F11").Value = 2621
$ws.Range("F14").Value = 216
$ws.Range("F15").Value = 6348
$ws.Range("F16").Value = 100
$ws.Range("F20").Value = 2549
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 6613
$ws.Range("F25").Value = 223
$ws.Range("F27").Value = 83
$ws.Range("F28").Value = 158
$ws.Range("F31").Value = 7267
$ws.Range("F34").Value = 45
$ws.Range("F43").Value = 16
$ws.Range("F46").Value = 3574
$ws.Range("F47").Value = 107
$ws.Range("F49").Value = 1142
$ws.Range("F51").Value = 42
